# Canpotex Data feed updated for 2023 Manual forecast
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update SEAsia forecast figures for rows 178-181 (Manual forecast refresh)
$ws.Range("B178").Value = 833
$ws.Range("B179").Value = 788
$ws.Range("B180").Value = 676
$ws.Range("B181").Value = 560

# Scroll the sheet view up and move the active selection to where the
# analyst was last working, matching the saved view state
$excel.ActiveWindow.ScrollRow = 166
$ws.Range("B170:B181").Select()
